# Commit message: "Added: Player Set target/limit FPS"
# This adds three new error-code rows to the ErrorCodes sheet:
#   - "Background colour not set"
#   - "Background colour wrong format"
#   - "Target FPS not set"
# The first two are inserted right before the existing "Trying to allocate
# to a full pool." row (pushing it down by two rows), and the third is
# appended immediately after that same row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 22:23 - this pushes the old row 22
# ("Trying to allocate to a full pool." / PoolAllocator / New & Alloc)
# down to row 24, keeping all of its existing cell values intact.
$ws.Rows("22:23").Insert()

# Insert one more blank row right after the (now shifted) pool-allocation
# row, at row 25.
$ws.Rows("25:25").Insert()

# New row 22: Background colour not set
$ws.Range("B22").Value = 20
$ws.Range("C22").Value = "Background colour not set"
$ws.Range("D22").Value = "Entry "
$ws.Range("E22").Value = "main"

# New row 23: Background colour wrong format
$ws.Range("B23").Value = 21
$ws.Range("C23").Value = "Background colour wrong format"
$ws.Range("D23").Value = "Entry "
$ws.Range("E23").Value = "main"

# Row 24 keeps its original text (pool allocator error) but the running
# error number needs to be bumped from 20 to 22 to stay sequential.
$ws.Range("B24").Value = 22

# New row 25: Target FPS not set
$ws.Range("B25").Value = 23
$ws.Range("C25").Value = "Target FPS not set"
$ws.Range("D25").Value = "Entry "
$ws.Range("E25").Value = "main"

# Match the author's final view/selection state.
$ws.Range("C26").Select() | Out-Null
